# Append 9 new rows (22-30) to the master-reg_center_user_machine_ sheet,
# matching the existing table's columns (A..H), and update the active
# selection to F14 as recorded by the saved workbook view state.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newRows = @(
    @(10002, 110021, 10021, "eng", $true, "superadmin", "now()", "now()"),
    @(10003, 110022, 10022, "eng", $true, "superadmin", "now()", "now()"),
    @(10004, 110023, 10023, "eng", $true, "superadmin", "now()", "now()"),
    @(10005, 110024, 10024, "eng", $true, "superadmin", "now()", "now()"),
    @(10006, 110025, 10025, "eng", $true, "superadmin", "now()", "now()"),
    @(10007, 110026, 10026, "eng", $true, "superadmin", "now()", "now()"),
    @(10008, 110027, 10027, "eng", $true, "superadmin", "now()", "now()"),
    @(10009, 110028, 10028, "eng", $true, "superadmin", "now()", "now()"),
    @(10010, 110029, 10029, "eng", $true, "superadmin", "now()", "now()")
)

$startRow = 22
for ($i = 0; $i -lt $newRows.Count; $i++) {
    $row = $startRow + $i
    $data = $newRows[$i]
    for ($c = 0; $c -lt $data.Count; $c++) {
        $ws.Cells.Item($row, $c + 1).Value = $data[$c]
    }
}

# Restore the saved selection/active cell.
$ws.Range("F14").Select() | Out-Null

# Set the print page orientation (records pageSetup settings on save).
$ws.PageSetup.Orientation = 1

$wb.Save() | Out-Null
